# Rename the two worksheets.
$wb = $excel.ActiveWorkbook
$wsAttr = $wb.Worksheets.Item(1)
$wsVals = $wb.Worksheets.Item(2)

$wsAttr.Name = "Attributes"
$wsVals.Name = "Possible values"

# --- Rebuild the "Possible values" sheet with the new lookup lists ---
$wsVals.Cells.Clear()

# Drop the old bestFit widths inherited from the previous layout (columns D
# and E used to hold the long "Filepath" example text) before laying the new
# columns out.
$wsVals.Columns.Item(4).ColumnWidth = $wsVals.StandardWidth
$wsVals.Columns.Item(5).ColumnWidth = $wsVals.StandardWidth

# Cells are written in the same order the original author typed them so the
# shared-string table is appended in the same sequence.
$wsVals.Range("A1").Value = "SPECTROMETER.Type"
$wsVals.Range("B1").Value = "VIPA"
$wsVals.Range("D1").Value = "TFP"
$wsVals.Range("E1").Value = "Time-Domain"
$wsVals.Range("F1").Value = "Stimulated"

$wsVals.Range("A2").Value = "SPECTROMETER.Detector_Type"
$wsVals.Range("B2").Value = "EMCCD"
$wsVals.Range("C2").Value = "CMOS"
$wsVals.Range("D2").Value = "sCMOS"
$wsVals.Range("E2").Value = "CCD"

$wsVals.Range("G1").Value = "uss-BM"
$wsVals.Range("C1").Value = "ar-BM"

$wsVals.Columns.Item(1).ColumnWidth = 32.5
$wsVals.Columns.Item(5).ColumnWidth = 11.6640625

$wsVals.Activate()
$wsVals.Range("G11").Select()

# --- Add data validation drop-downs on the "Attributes" sheet ---
$wsAttr.Range("B14").Validation.Delete()
$wsAttr.Range("B14").Validation.Add(3, 1, 1, "='Possible values'!`$B`$1:`$Z`$1")
$wsAttr.Range("B14").Validation.IgnoreBlank = $true
$wsAttr.Range("B14").Validation.InCellDropdown = $true
$wsAttr.Range("B14").Validation.ShowInput = $true
$wsAttr.Range("B14").Validation.ShowError = $true

$wsAttr.Range("B20").Validation.Delete()
$wsAttr.Range("B20").Validation.Add(3, 1, 1, "='Possible values'!`$B`$2:`$Z`$2")
$wsAttr.Range("B20").Validation.IgnoreBlank = $true
$wsAttr.Range("B20").Validation.InCellDropdown = $true
$wsAttr.Range("B20").Validation.ShowInput = $true
$wsAttr.Range("B20").Validation.ShowError = $true

# --- Restore the view / selection on the "Attributes" sheet ---
$wsAttr.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$wsAttr.Range("B21").Select()
